$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 399
$ws.Range("I38").Value = 399
$ws.Range("K38").Value = 1197
$ws.Range("M38").Value = -825
$ws.Range("H55").Value = 154
$ws.Range("I55").Value = 169.75
$ws.Range("K55").Value = 169.75
$ws.Range("M55").Value = 44.25
$ws.Range("H86").Value = 14050.3
$ws.Range("I86").Value = 14050.3
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 14050.3
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -12927.3
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 14050.3
$ws.Range("I89").Value = 14050.3
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 70251.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -64635.5
$ws.Range("N89").ClearContents()
$ws.Range("H101").Value = 3281.6667
$ws.Range("I101").Value = 3518
$ws.Range("J101").Value = 2100
$ws.Range("K101").Value = 10554
$ws.Range("L101").Value = 6300
$ws.Range("M101").Value = -8932
$ws.Range("N101").Value = -9544
$ws.Range("H113").Value = 71429770
$ws.Range("I113").Value = 33334266
$ws.Range("J113").Value = 100001390
$ws.Range("K113").Value = 33334266
$ws.Range("L113").Value = 100001390
$ws.Range("M113").Value = -33331012
$ws.Range("N113").Value = -100007898
$ws.Range("H138").Value = 2256.7058
$ws.Range("I138").Value = 1352.8334
$ws.Range("J138").Value = 2749.7273
$ws.Range("K138").Value = 4058.5002
$ws.Range("L138").Value = 8249.1819
$ws.Range("M138").Value = 1081.4998
$ws.Range("N138").Value = -18529.1819
$ws.Range("H141").Value = 3183.3333
$ws.Range("I141").Value = 2773.75
$ws.Range("K141").Value = 8321.25
$ws.Range("M141").Value = -3141.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13890776
$ws.Range("I32").Value = 14287621
$ws.Range("K32").Value = 14287621
$ws.Range("M32").Value = -14287334
$ws.Range("H61").Value = 26375330
$ws.Range("I61").Value = 38466450
$ws.Range("J61").Value = 177909.67
$ws.Range("K61").Value = 38466450
$ws.Range("L61").Value = 177909.67
$ws.Range("M61").Value = -38466238
$ws.Range("N61").Value = -178333.67
$ws.Range("H97").Value = 1313
$ws.Range("I97").Value = 1168.4286
$ws.Range("J97").Value = 2325
$ws.Range("K97").Value = 1168.4286
$ws.Range("L97").Value = 2325
$ws.Range("M97").Value = -672.4286
$ws.Range("N97").Value = -3317
$ws.Range("H102").Value = 13453.962
$ws.Range("I102").Value = 19600.334
$ws.Range("K102").Value = 19600.334
$ws.Range("M102").Value = -17978.334
$ws.Range("H132").Value = 6254.2856
$ws.Range("I132").Value = 3836.0527
$ws.Range("K132").Value = 11508.1581
$ws.Range("M132").Value = -8978.158100000001
$ws.Range("H135").Value = 400000
$ws.Range("J135").Value = 400000
$ws.Range("L135").Value = 400000
$ws.Range("N135").Value = -410140
$ws.Range("H136").Value = 26375330
$ws.Range("I136").Value = 38466450
$ws.Range("J136").Value = 177909.67
$ws.Range("K136").Value = 115399350
$ws.Range("L136").Value = 533729.01
$ws.Range("M136").Value = -115396800
$ws.Range("N136").Value = -538829.01

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 9999.666999999999
$ws.Range("J22").Value = 4999
$ws.Range("L22").Value = 4999
$ws.Range("N22").Value = -5345
$ws.Range("H63").Value = 65135
$ws.Range("J63").Value = 65135
$ws.Range("L63").Value = 65135
$ws.Range("N63").Value = -66507
$ws.Range("H66").Value = 65135
$ws.Range("J66").Value = 65135
$ws.Range("L66").Value = 195405
$ws.Range("N66").Value = -202269
$ws.Range("H94").Value = 1810.0435
$ws.Range("I94").Value = 1814.591
$ws.Range("K94").Value = 1814.591
$ws.Range("M94").Value = -1363.591

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2231.1667
$ws.Range("I16").Value = 2189.2856
$ws.Range("K16").Value = 2189.2856
$ws.Range("M16").Value = -1902.2856
$ws.Range("H22").Value = 329.6
$ws.Range("I22").Value = 329.6
$ws.Range("K22").Value = 329.6
$ws.Range("M22").Value = 20.39999999999998
$ws.Range("H98").Value = 55963
$ws.Range("J98").Value = 55963
$ws.Range("L98").Value = 55963
$ws.Range("N98").Value = -60455
$ws.Range("H112").Value = 75506.336
$ws.Range("J112").Value = 75506.336
$ws.Range("L112").Value = 75506.336
$ws.Range("N112").Value = -78460.336
$ws.Range("H113").Value = 2231.1667
$ws.Range("I113").Value = 2189.2856
$ws.Range("K113").Value = 2189.2856
$ws.Range("M113").Value = -19.28560000000016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6743.6924
$ws.Range("I131").Value = 5978.909
$ws.Range("J131").Value = 10950
$ws.Range("K131").Value = 17936.727
$ws.Range("L131").Value = 32850
$ws.Range("M131").Value = -12896.727
$ws.Range("N131").Value = -42930

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2826.125
$ws.Range("I80").Value = 3037.5454
$ws.Range("J80").Value = 2361
$ws.Range("K80").Value = 3037.5454
$ws.Range("L80").Value = 2361
$ws.Range("M80").Value = -2039.5454
$ws.Range("N80").Value = -4357
$ws.Range("H83").Value = 2826.125
$ws.Range("I83").Value = 3037.5454
$ws.Range("J83").Value = 2361
$ws.Range("K83").Value = 15187.727
$ws.Range("L83").Value = 11805
$ws.Range("M83").Value = -10195.727
$ws.Range("N83").Value = -21789
$ws.Range("H107").Value = 2109.2727
$ws.Range("I107").Value = 1299.5555
$ws.Range("J107").Value = 5753
$ws.Range("K107").Value = 1299.5555
$ws.Range("L107").Value = 5753
$ws.Range("M107").Value = 620.4445000000001
$ws.Range("N107").Value = -9593
$ws.Range("H113").Value = 3225.1304
$ws.Range("I113").Value = 3378.1765
$ws.Range("J113").Value = 2791.5
$ws.Range("K113").Value = 3378.1765
$ws.Range("L113").Value = 2791.5
$ws.Range("M113").Value = -1208.1765
$ws.Range("N113").Value = -7131.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2916.7144
$ws.Range("I22").Value = 2964.3572
$ws.Range("K22").Value = 2964.3572
$ws.Range("M22").Value = -2669.3572
$ws.Range("H27").Value = 2916.7144
$ws.Range("I27").Value = 2964.3572
$ws.Range("K27").Value = 2964.3572
$ws.Range("M27").Value = -2857.3572
$ws.Range("H46").Value = 4325.381
$ws.Range("I46").Value = 3967.923
$ws.Range("J46").Value = 4906.25
$ws.Range("K46").Value = 3967.923
$ws.Range("L46").Value = 4906.25
$ws.Range("M46").Value = -3779.923
$ws.Range("N46").Value = -5282.25
$ws.Range("H55").Value = 62500264
$ws.Range("I55").Value = 83333600
$ws.Range("J55").Value = 258.25
$ws.Range("K55").Value = 83333600
$ws.Range("L55").Value = 258.25
$ws.Range("M55").Value = -83333427
$ws.Range("N55").Value = -604.25
$ws.Range("I61").Value = 1167.3334
$ws.Range("J61").Value = 1001
$ws.Range("K61").Value = 1167.3334
$ws.Range("L61").Value = 1001
$ws.Range("M61").Value = -965.3334
$ws.Range("N61").Value = -1405
$ws.Range("I113").Value = 1167.3334
$ws.Range("J113").Value = 1001
$ws.Range("K113").Value = 1167.3334
$ws.Range("L113").Value = 1001
$ws.Range("M113").Value = 1002.6666
$ws.Range("N113").Value = -5341
$ws.Range("H130").Value = 24997
$ws.Range("J130").Value = 24997
$ws.Range("L130").Value = 24997
$ws.Range("N130").Value = -35037

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H122").Value = 4061.2144
$ws.Range("I122").Value = 1824.8077
$ws.Range("K122").Value = 5474.4231
$ws.Range("M122").Value = -3024.4231
$ws.Range("H130").Value = 88495
$ws.Range("J130").Value = 88495
$ws.Range("L130").Value = 88495
$ws.Range("N130").Value = -98535
